# Adds a new "2022" column (H) to the sheet, mirroring the formatting of the
# existing "2021" column (G), row by row, and populates the new figures.
# Also updates the active selection to I4 (the cell just past the new data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new H-column value ($null for rows that only carry formatting,
# i.e. the section-header rows which are blank in every year column).
$values = @{
    4  = 2022
    5  = 92.960099223795225
    6  = $null
    7  = 96.03949422949897
    8  = 91.012153547624152
    9  = $null
    10 = 94.391087218067838
    11 = 91.76755842559642
    12 = $null
    13 = 92.942689638142156
    14 = 86.897877953385489
    15 = 96.500794494289821
    16 = 94.135975315309977
    17 = 89.456106196597958
    18 = 94.270923428904894
    19 = 97.027480110114013
    20 = 98.077227596867303
    21 = 90.983384827072243
    22 = $null
    23 = 90.468970496790078
    24 = 95.809965597614095
    25 = 88.221110530662017
    26 = $null
    27 = 69.811292606515579
    28 = 85.757158930558518
    29 = 93.032103866435918
    30 = 97.325262246493097
    31 = 98.908492141713779
    32 = $null
    33 = 91.968006037496949
    34 = 91.809335747904541
    35 = 91.27524653351901
    36 = 94.397288657466234
    37 = 96.740699993405215
}

foreach ($row in 4..37) {
    # Copy column G's formatting (number format, font, borders, alignment)
    # for this row onto the new column H cell, then write its value (if any).
    $ws.Range("G$row").Copy()
    $ws.Range("H$row").PasteSpecial(-4122)

    $val = $values[$row]
    if ($val -ne $null) {
        $ws.Range("H$row").Value = $val
    }
}

$excel.CutCopyMode = $false

# Move the active selection to I4, matching the saved view state.
$ws.Range("I4").Select() | Out-Null
